$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "2020-12-16 00:00:00"
$ws.Range("I2").Value = 0.08
$ws.Range("P2").Value = 1.068917763333
$ws.Range("Q2").Value = 0.174051527917
